$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) columns
# D-column cells get NumberFormat "@" (text) first so Excel does not
# reinterpret dotted/zero-padded price strings (e.g. "0.5600", "1.000") as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.086.42"
$ws.Range("E2").Value = "  +0.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.890.27"
$ws.Range("E3").Value = "  +1.43%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.92"
$ws.Range("E5").Value = "  +0.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5154"
$ws.Range("E7").Value = "  +1.77%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3756"
$ws.Range("E8").Value = "  +2.96%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07200"
$ws.Range("E9").Value = "  +0.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.12"
$ws.Range("E10").Value = "  +1.32%  "

$ws.Range("E11").Value = "  +0.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07662"
$ws.Range("E12").Value = "  +2.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.889.87"
$ws.Range("E13").Value = "  +1.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.32"
$ws.Range("E14").Value = "  +1.81%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.248"
$ws.Range("E15").Value = "  +0.18%  "

$ws.Range("E16").Value = "  +0.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008505"
$ws.Range("E17").Value = "  +0.03%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.46"
$ws.Range("E18").Value = "  +1.84%  "

$ws.Range("E19").Value = "  +0.12%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.128.03"
$ws.Range("E20").Value = "  +0.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.062"
$ws.Range("E21").Value = "  +0.41%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.123.19"
$ws.Range("E22").Value = "  +3.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.58"
$ws.Range("E23").Value = "  +1.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.381"
$ws.Range("E24").Value = "  -0.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.313"
$ws.Range("E25").Value = "  +11.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.96"
$ws.Range("E26").Value = "  -0.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.765"
$ws.Range("E27").Value = "  -1.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.07"
$ws.Range("E28").Value = "  +1.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.32"
$ws.Range("E29").Value = "  +1.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.931"
$ws.Range("E30").Value = "  +5.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.798"
$ws.Range("E31").Value = "  +2.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09203"
$ws.Range("E32").Value = "  -0.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05057"
$ws.Range("E33").Value = "  -1.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.235"
$ws.Range("E34").Value = "  +7.15%  "

$ws.Range("E35").Value = "  +2.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.980"
$ws.Range("E36").Value = "  +0.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.289"
$ws.Range("E37").Value = "  +0.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.603"
$ws.Range("E38").Value = "  +2.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5600"
$ws.Range("E39").Value = "  +1.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01994"
$ws.Range("E40").Value = "  -0.42%  "

$ws.Range("E41").Value = "  +0.36%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.058"
$ws.Range("E42").Value = "  +6.21%  "

$ws.Range("E43").Value = "  +1.91%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "118.19"
$ws.Range("E44").Value = "  -0.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1504"
$ws.Range("E45").Value = "  +2.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4817"
$ws.Range("E46").Value = "  +2.66%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.17"
$ws.Range("E47").Value = "  +0.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.000"
$ws.Range("E48").Value = "  +0.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.603"
$ws.Range("E49").Value = "  +2.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.61"
$ws.Range("E50").Value = "  +2.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.09"
$ws.Range("E51").Value = "  +1.84%  "
